$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add "Country" label in A3
$ws.Range("A3").Value = "Country"

# Fill A4:A45 with "Belgium" for every data row (establishes shared string
# index for "Belgium" before "Blank Node" is introduced)
for ($r = 4; $r -le 45; $r++) {
    $ws.Cells.Item($r, 1).Value = "Belgium"
}

# Blank Node placeholder header in D3
$ws.Range("D3").Value = "Blank Node"

# Restore the active selection to D3, as recorded by the original author
$ws.Range("D3").Select()
